# edit.ps1 - applies the "Small updates - transition to marker if not on
# visible map" revision to "Specifications from Udacity.docx".
#
# Summary of the changes applied below:
#  1. Move the hidden "_GoBack" bookmark from around the word "map" (in the
#     "For sake of efficiency, the map API..." bullet) to wrap
#     "Neighborhood Map " in the very first bulleted item.
#  2. Merge several runs that used to be split into multiple <w:r> elements
#     back into single runs (pure whitespace/run-boundary cleanup, no text
#     changes): the CRITERIA heading, Responsiveness / Usability bullets,
#     Filter Locations bullet, "Map displays..." bullet, "There are at
#     least 5 locations..." bullet, README / Comments bullets and the
#     JavaScript Style Guide bullet.
#  3. Register a (previously only latent) "Emphasis" character style in
#     styles.xml, matching Word's built-in definition.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Bookmark relocation
# ---------------------------------------------------------------------

# 1a. Delete the existing hidden bookmark (currently wraps "map" in the
#     "For sake of efficiency, the map API..." bullet).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# 1b. Re-create it around "Neighborhood Map " in the first bulleted item.
$introPara = $d.Paragraphs(2)
$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End)
$introRange.Find.Execute("Neighborhood Map ") | Out-Null
$d.Bookmarks.Add("_GoBack", $introRange) | Out-Null

# ---------------------------------------------------------------------
# 2. Run-merge clean ups (Find/Replace with identical text re-flows the
#    run boundaries without altering the visible text).
# ---------------------------------------------------------------------

function Merge-Run([int]$paraIndex, [string]$text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# CRITERIA heading: "CRITERIA" + " " -> "CRITERIA "; "for  " + "MEETS" -> "for  MEETS"
Merge-Run 48 "CRITERIA "
Merge-Run 48 "for  MEETS"

# Responsiveness / Usability bullets
Merge-Run 50 "Responsiveness - All application components render on-screen in a responsive manner."
Merge-Run 51 "Usability - All application components are usable across modern desktop, tablet, and phone browsers."

# Filter Locations bullet (also absorbs the "markers" spell-check run)
Merge-Run 54 " - Includes a text input field or dropdown menu that filters the map markers and list items to locations matching the text input or selection. Filter function runs error-free."

# "Map displays..." bullet
Merge-Run 60 "Map displays all location markers by default, and displays the filtered subset of location markers when a filter is applied."

# "There are at least 5 locations..." bullet
Merge-Run 67 "There are at least 5 locations in the model. These may be hard-coded or retrieved from a data API"

# README / Comments bullets
Merge-Run 84 " - A README file is included detailing all steps required to successfully run the application."
Merge-Run 85 " - Comments are present and effectively explain longer code procedures."

# JavaScript Style Guide bullet (merge trailing space run)
Merge-Run 87 " JavaScript Style Guide. "

# ---------------------------------------------------------------------
# 3. Register the "Emphasis" character style.
#
# Word only writes a built-in style's full definition into styles.xml once
# it is actually used somewhere in the document. To materialize it with
# the exact same shape Word itself produces (type="character",
# basedOn="DefaultParagraphFont", uiPriority 20, qFormat, italic +
# italic-complex-script), we briefly apply it to a throwaway character and
# then restore that paragraph's original text/run exactly as it was.
# ---------------------------------------------------------------------

$scratchPara = $d.Paragraphs(1)
$scratchFull = $scratchPara.Range
$originalText = $d.Range($scratchFull.Start, $scratchFull.End - 1).Text

$scratchChar = $d.Range($scratchFull.Start, $scratchFull.Start + 1)
$scratchChar.Style = "Emphasis"

$emphasis = $d.Styles("Emphasis")
$emphasis.Font.Italic = $true
$emphasis.Font.ItalicBi = $true
$emphasis.Priority = 20

# Restore paragraph 1's text/run exactly (strip the temporary rStyle).
$restoreRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.End - 1)
$restoreRange.Delete()
$restoreRange.InsertBefore($originalText)

Write-Output "done"
